# Update the "Förändrad" (Changed) date column (C) for rows 2-10
# from serial 45184 (2023-09-15) to serial 45185 (2023-09-16),
# keeping existing number formatting / style intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45184) {
        $cell.Value2 = 45185
    }
}
